# Edit GPCRmd_B2AR_nomenclature_test.xlsx:
#  - Fill in the previously-blank rows 6-14 with TM1 nomenclature data
#  - Move the "new TM segment" marker row (single bold cell) up by one row,
#    for each of the 7 segment boundaries later in the sheet
#  - Remove the now-unused trailing row 256 (sheet shrinks from 256 to 255 rows)
#  - Select cell A15, as left by the editor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Populate rows 6-14 (continuation of the TM1 nomenclature table)
# ---------------------------------------------------------------------------
$tm1Rows = @(
    @("1x29", 1.39, "M40"),
    @("1x30", 1.40, "S41"),
    @("1x31", 1.41, "L42"),
    @("1x32", 1.42, "I43"),
    @("1x33", 1.43, "V44"),
    @("1x34", 1.44, "L45"),
    @("1x35", 1.45, "A46"),
    @("1x36", 1.46, "I47"),
    @("1x37", 1.47, "V48")
)

$r = 6
foreach ($row in $tm1Rows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    # Match the row height used by the rest of the nomenclature table (ht=14)
    $ws.Rows.Item($r).RowHeight = 14
    $r++
}

# ---------------------------------------------------------------------------
# 2) Shift the TM-segment marker rows up by one
#    (before: row X-1 = plain s=2 triplet, row X = lone bold s=1 marker)
#    (after : row X-1 = lone bold s=1 marker, row X = plain s=2 triplet)
# ---------------------------------------------------------------------------
$markerRows = @(38, 70, 107, 135, 178, 217, 243)

foreach ($mr in $markerRows) {
    $prev = $mr - 1

    # Row (prev) currently holds the plain s=2 triplet formatting; stamp
    # that same formatting onto row (mr) first, while it is still available.
    $ws.Range("A$prev`:C$prev").Copy() | Out-Null
    $ws.Range("A$mr`:C$mr").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

    # Now turn cell A(prev) into the bold lone marker cell (copy formatting
    # from the known s=1 styled header cell A1) and drop its B/C neighbours.
    $ws.Range("A1").Copy() | Out-Null
    $ws.Range("A$prev").PasteSpecial(-4122) | Out-Null           # xlPasteFormats
    $ws.Range("B$prev`:C$prev").Clear() | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Drop the now-superfluous last row (256)
# ---------------------------------------------------------------------------
$ws.Rows.Item(256).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4) Leave the selection on A15, matching the saved view state
# ---------------------------------------------------------------------------
$ws.Range("A15").Select() | Out-Null
